$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns("C").Insert()
$ws.Range("C2").Value = "Value"
$ws.Range("C5").Value = "status"
$ws.Columns("C").AutoFit() | Out-Null
Write-Host "done"
